$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 121.7744789210046
$ws.Range("C2").Value = 40
$ws.Range("D2").Value = 212.815665118663
$ws.Range("F2").Value = 2102.493658644074
$ws.Range("H2").Value = 610.5473328328467
$ws.Range("J2").Value = 0.002348093274886724
$ws.Range("K2").Value = 0.05205780726319665
$ws.Range("L2").Value = 0.01936550632414199
$ws.Range("N2").Value = 1011.226663093531
$ws.Range("P2").Value = 0.93978044280003
$ws.Range("S2").Value = 0.8275679403271208
$ws.Range("T2").Value = 1739.956346634909
$ws.Range("V2").Value = 0.4429960375016347
$ws.Range("W2").Value = 1011.226663093531
$ws.Range("Y2").Value = 0.5613696935541674
$ws.Range("Z2").Value = 342.7427691326894
$ws.Range("B3").Value = 210.6375000841791
$ws.Range("C3").Value = 40
$ws.Range("D3").Value = 216.8082181648471
$ws.Range("F3").Value = 2406.135090620478
$ws.Range("H3").Value = 603.2409607583298
$ws.Range("J3").Value = 0.02220982661171617
$ws.Range("K3").Value = 0.08055904135278205
$ws.Range("L3").Value = 0.01986982261615612
$ws.Range("N3").Value = 994.5797031875584
$ws.Range("P3").Value = 0.9224742571701446
$ws.Range("S3").Value = 0.8330705843437558
$ws.Range("T3").Value = 2004.480365953217
$ws.Range("V3").Value = 0.3932621536389742
$ws.Range("W3").Value = 994.5797031875584
$ws.Range("Y3").Value = 0.5380433214428376
$ws.Range("Z3").Value = 324.5697701567802
$ws.Range("B4").Value = 99.04886500521728
$ws.Range("D4").Value = 216.9987957756236
$ws.Range("F4").Value = 1941.622773228155
$ws.Range("H4").Value = 602.8922037306088
$ws.Range("J4").Value = 0.0113203468964704
$ws.Range("K4").Value = 0.04469268339914839
$ws.Range("L4").Value = 0.01628156976606811
$ws.Range("N4").Value = 1008.049311634372
$ws.Range("P4").Value = 0.9216641008773141
$ws.Range("S4").Value = 0.8433294465746841
$ws.Range("T4").Value = 1637.427658803304
$ws.Range("V4").Value = 0.4798330029566718
$ws.Range("W4").Value = 1008.049311634372
$ws.Range("Y4").Value = 0.4369039344566191
$ws.Range("Z4").Value = 263.4059758631246
$ws.Range("B5").Value = 94.44367166706512
$ws.Range("C5").Value = 16.09663463598934
$ws.Range("D5").Value = 219.8993802730784
$ws.Range("F5").Value = 2853.159315040723
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0.008552758480979165
$ws.Range("K5").Value = 0.03244225616477758
$ws.Range("L5").Value = 0.03804039643536632
$ws.Range("N5").Value = 1013.306053125703
$ws.Range("P5").Value = 0.9095068833374305
$ws.Range("S5").Value = 0.840365493926777
$ws.Range("T5").Value = 2397.696637035982
$ws.Range("V5").Value = 0.3900786714964157
$ws.Range("W5").Value = 1013.306053125703
$ws.Range("Y5").Value = 0.5162579709204294
$ws.Range("Z5").Value = 0
$ws.Range("B6").Value = 131.3991036849226
$ws.Range("C6").Value = 40
$ws.Range("D6").Value = 216.4676429633167
$ws.Range("F6").Value = 2482.269757635597
$ws.Range("H6").Value = 603.8642133771305
$ws.Range("J6").Value = 0.02171922704753687
$ws.Range("K6").Value = 0.04876516460272357
$ws.Range("L6").Value = 0.009355249664325332
$ws.Range("N6").Value = 988.7071353895649
$ws.Range("P6").Value = 0.9239256142955861
$ws.Range("S6").Value = 0.8337761221893575
$ws.Range("T6").Value = 2069.657252749324
$ws.Range("V6").Value = 0.3904080811222203
$ws.Range("W6").Value = 988.7071353895649
$ws.Range("Y6").Value = 0.4352623232209656
$ws.Range("Z6").Value = 262.8393404245307
$ws.Range("B7").Value = 132.6946672458907
$ws.Range("C7").Value = 40.00000000000001
$ws.Range("D7").Value = 0
$ws.Range("F7").Value = 2148.750086759504
$ws.Range("H7").Value = 1000
$ws.Range("J7").Value = 0.01038766567243881
$ws.Range("K7").Value = 0.05864246512268719
$ws.Range("L7").Value = 0.006686665225725419
$ws.Range("N7").Value = 1015.651705141702
$ws.Range("P7").Value = 0.9414457470998366
$ws.Range("Q7").Value = 0
$ws.Range("S7").Value = 0.8287895735837484
$ws.Range("T7").Value = 1780.861668143451
$ws.Range("V7").Value = 0.4385542201385005
$ws.Range("W7").Value = 1015.651705141702
$ws.Range("Y7").Value = 0.5350474224610067
$ws.Range("Z7").Value = 535.0474224610067
$ws.Range("B8").Value = 84.26707576041001
$ws.Range("C8").Value = 40
$ws.Range("D8").Value = 215.0116188266796
$ws.Range("F8").Value = 1932.703289267768
$ws.Range("H8").Value = 606.5287375471762
$ws.Range("J8").Value = 0.02160526545129312
$ws.Range("K8").Value = 0.03613067964798337
$ws.Range("L8").Value = 0.01614401537016744
$ws.Range("N8").Value = 990.6938417686193
$ws.Range("P8").Value = 0.930182290107864
$ws.Range("S8").Value = 0.8602370220641031
$ws.Range("T8").Value = 1662.582922093202
$ws.Range("V8").Value = 0.451632929765607
$ws.Range("W8").Value = 990.6938417686193
$ws.Range("Y8").Value = 0.5457268637451181
$ws.Range("Z8").Value = 330.9990257129064
$ws.Range("B9").Value = 95.81084823624377
$ws.Range("D9").Value = 219.5232869633142
$ws.Range("F9").Value = 2235.171936789313
$ws.Range("H9").Value = 598.2723848571351
$ws.Range("J9").Value = 0.00310669693234989
$ws.Range("K9").Value = 0.04007753794520412
$ws.Range("L9").Value = 0.009274484237896724
$ws.Range("N9").Value = 989.6338231384073
$ws.Range("P9").Value = 0.911065075448798
$ws.Range("Q9").Value = 200
$ws.Range("S9").Value = 0.8480620151444597
$ws.Range("T9").Value = 1895.56441690789
$ws.Range("V9").Value = 0.4102260702268819
$ws.Range("W9").Value = 989.6338231384073
$ws.Range("Y9").Value = 0.5296023049271085
$ws.Range("Z9").Value = 316.8464339945768
$ws.Range("B10").Value = 165.5847935824491
$ws.Range("C10").Value = 40
$ws.Range("D10").Value = 217.9593461984611
$ws.Range("F10").Value = 2550.877625225707
$ws.Range("H10").Value = 601.1343964568161
$ws.Range("J10").Value = 0.03112316767652012
$ws.Range("K10").Value = 0.05450283808351926
$ws.Range("L10").Value = 0.03288971358490739
$ws.Range("N10").Value = 983.6246111427911
$ws.Range("P10").Value = 0.9176023120288297
$ws.Range("S10").Value = 0.843787696727783
$ws.Range("T10").Value = 2152.399156023636
$ws.Range("V10").Value = 0.3645456662837667
$ws.Range("W10").Value = 983.6246111427911
$ws.Range("Y10").Value = 0.5752813416454232
$ws.Range("Z10").Value = 345.821402102889
$ws.Range("B11").Value = 87.87486582850033
$ws.Range("D11").Value = 224.2148775141743
$ws.Range("F11").Value = 2640.103597130349
$ws.Range("H11").Value = 589.6867741490611
$ws.Range("J11").Value = 0.02329551994717612
$ws.Range("K11").Value = 0.0286362862600064
$ws.Range("L11").Value = 0.01195363644299093
$ws.Range("N11").Value = 983.9256373026786
$ws.Range("P11").Value = 0.8920014684902277
$ws.Range("Q11").Value = 200
$ws.Range("S11").Value = 0.8449113856185803
$ws.Range("T11").Value = 2230.653588428001
$ws.Range("V11").Value = 0.3644502356478774
$ws.Range("W11").Value = 983.9256373026786
$ws.Range("Y11").Value = 0.4563430603566268
$ws.Range("Z11").Value = 269.0994671670095
